$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 78 with Q3 2021 data (update from MV source data)
$ws.Range("A78").Value = "'01-07-2021"
$ws.Range("A78").Style = "Normal"
$ws.Range("C78").Value = -0.8
$ws.Range("D78").Value = 0.89
$ws.Range("E78").Value = 1.7
$ws.Range("F78").Value = 2.32
$ws.Range("G78").Value = 2.46
